$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1849148418491484
$ws.Range("C2").Value = 0.5815085158150851
$ws.Range("P2").Value = 0.1289537712895377
$ws.Range("S2").Value = 0.1046228710462287
$ws.Range("B3").Value = 0.01673640167364017
$ws.Range("C3").Value = 0.01673640167364017
$ws.Range("J3").Value = 0.008368200836820083
$ws.Range("P3").Value = 0.7405857740585774
$ws.Range("S3").Value = 0.2175732217573222
$ws.Range("J4").Value = 0.03703703703703703
$ws.Range("P4").Value = 0.6481481481481481
$ws.Range("S4").Value = 0.3148148148148148
$ws.Range("B6").Value = 0.06909090909090909
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("E6").Value = 0.007272727272727273
$ws.Range("F6").Value = 0.04
$ws.Range("J6").Value = 0.3345454545454545
$ws.Range("O6").Value = 0.02909090909090909
$ws.Range("Q6").Value = 0.1054545454545455
$ws.Range("R6").Value = 0.05454545454545454
$ws.Range("S6").Value = 0.3418181818181818
$ws.Range("B7").Value = 0.1400778210116732
$ws.Range("D7").Value = 0.01167315175097276
$ws.Range("F7").Value = 0.02723735408560311
$ws.Range("J7").Value = 0.1439688715953307
$ws.Range("O7").Value = 0.0311284046692607
$ws.Range("Q7").Value = 0.132295719844358
$ws.Range("R7").Value = 0.07782101167315175
$ws.Range("S7").Value = 0.4357976653696498
$ws.Range("B8").Value = 0.1090909090909091
$ws.Range("D8").Value = 0.02479338842975207
$ws.Range("F8").Value = 0.05785123966942149
$ws.Range("J8").Value = 0.112396694214876
$ws.Range("O8").Value = 0.02479338842975207
$ws.Range("Q8").Value = 0.1900826446280992
$ws.Range("R8").Value = 0.06776859504132231
$ws.Range("S8").Value = 0.4132231404958678
$ws.Range("B9").Value = 0.09183673469387756
$ws.Range("D9").Value = 0.02040816326530612
$ws.Range("E9").Value = 0.00510204081632653
$ws.Range("F9").Value = 0.05102040816326531
$ws.Range("J9").Value = 0.163265306122449
$ws.Range("O9").Value = 0.01020408163265306
$ws.Range("Q9").Value = 0.1785714285714286
$ws.Range("R9").Value = 0.04591836734693878
$ws.Range("S9").Value = 0.4336734693877551
$ws.Range("B10").Value = 0.1294765840220386
$ws.Range("D10").Value = 0.01928374655647383
$ws.Range("E10").Value = 0.0006887052341597796
$ws.Range("F10").Value = 0.06818181818181818
$ws.Range("J10").Value = 0.1260330578512397
$ws.Range("O10").Value = 0.01377410468319559
$ws.Range("Q10").Value = 0.2403581267217631
$ws.Range("R10").Value = 0.05234159779614325
$ws.Range("S10").Value = 0.349862258953168
$ws.Range("G11").Value = 0.1204188481675393
$ws.Range("J11").Value = 0.1178010471204188
$ws.Range("K11").Value = 0.1858638743455497
$ws.Range("L11").Value = 0.5549738219895288
$ws.Range("S11").Value = 0.02094240837696335
$ws.Range("G12").Value = 0.7822222222222223
$ws.Range("J12").Value = 0.1422222222222222
$ws.Range("K12").Value = 0.004444444444444444
$ws.Range("L12").Value = 0.02222222222222222
$ws.Range("S12").Value = 0.04888888888888889
$ws.Range("F13").Value = 0.01449275362318841
$ws.Range("G13").Value = 0.6231884057971014
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.1014492753623188
$ws.Range("F15").Value = 0.02788844621513944
$ws.Range("H15").Value = 0.1792828685258964
$ws.Range("I15").Value = 0.04780876494023904
$ws.Range("J15").Value = 0.2788844621513944
$ws.Range("K15").Value = 0.05976095617529881
$ws.Range("M15").Value = 0.01195219123505976
$ws.Range("N15").Value = 0.00796812749003984
$ws.Range("O15").Value = 0.0796812749003984
$ws.Range("S15").Value = 0.3067729083665339
$ws.Range("F16").Value = 0.03422053231939164
$ws.Range("H16").Value = 0.1673003802281369
$ws.Range("I16").Value = 0.07224334600760456
$ws.Range("J16").Value = 0.3422053231939163
$ws.Range("K16").Value = 0.1368821292775665
$ws.Range("M16").Value = 0.01520912547528517
$ws.Range("O16").Value = 0.05703422053231939
$ws.Range("S16").Value = 0.1749049429657795
$ws.Range("F17").Value = 0.04225352112676056
$ws.Range("H17").Value = 0.2007042253521127
$ws.Range("I17").Value = 0.08626760563380281
$ws.Range("J17").Value = 0.3485915492957746
$ws.Range("K17").Value = 0.09330985915492958
$ws.Range("M17").Value = 0.02112676056338028
$ws.Range("N17").Value = 0.00176056338028169
$ws.Range("O17").Value = 0.05633802816901409
$ws.Range("S17").Value = 0.1496478873239437
$ws.Range("F18").Value = 0.0245398773006135
$ws.Range("H18").Value = 0.2085889570552147
$ws.Range("I18").Value = 0.0736196319018405
$ws.Range("J18").Value = 0.392638036809816
$ws.Range("K18").Value = 0.0736196319018405
$ws.Range("M18").Value = 0.006134969325153374
$ws.Range("O18").Value = 0.1042944785276074
$ws.Range("S18").Value = 0.1165644171779141
$ws.Range("F19").Value = 0.02626641651031895
$ws.Range("H19").Value = 0.2338961851156973
$ws.Range("I19").Value = 0.06504065040650407
$ws.Range("J19").Value = 0.3352095059412132
$ws.Range("K19").Value = 0.1169480925578487
$ws.Range("M19").Value = 0.03064415259537211
$ws.Range("O19").Value = 0.05128205128205128
$ws.Range("S19").Value = 0.1407129455909944
